$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13 and 14: coin data swapped (WrappedEther <-> Litecoin) with updated price/volume
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'96.55"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.901.55"
$ws.Range("E14").Value = "  +0.13%  "

# Rows with both Price (D) and Volume (E) changes
$ws.Range("D2").Value = "27.225.37"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.894.76"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'307.78"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.5201"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'0.3768"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").Value = "'0.07277"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").Value = "'21.19"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "'0.9006"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "'0.08188"
$ws.Range("E12").Value = "  +6.92%  "
$ws.Range("D15").Value = "'5.288"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'0.000008597"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "'14.55"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D20").Value = "27.248.76"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").Value = "'5.085"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D23").Value = "'6.392"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "'2.305"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "'147.39"
$ws.Range("D27").Value = "'1.744"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "'115.29"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "'4.820"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").Value = "'4.914"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D32").Value = "'0.7966"
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").Value = "'0.05031"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "'1.220"
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").Value = "'3.447"
$ws.Range("E35").Value = "  +4.95%  "
$ws.Range("D36").Value = "'2.955"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "'2.590"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.5698"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").Value = "'0.01989"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'1.075"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'8.950"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'6.561"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "'115.56"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").Value = "'0.1514"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'0.4873"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D47").Value = "'10.05"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").Value = "'1.619"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").Value = "'38.19"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'63.44"
$ws.Range("E50").Value = "  -0.75%  "

# Rows with only Volume (E) changes
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("E51").Value = "  +0.34%  "
